$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Fill in "完成情况" (C column) for the existing last block (rows 164-169)
# ------------------------------------------------------------------
$ws.Range("C164").Value = "未完成"
$ws.Range("C165").Value = "未完成"
$ws.Range("C166").Value = "完成"
$ws.Range("C167").Value = "未完成"
$ws.Range("C168").Value = "未完成"
$ws.Range("C169").Value = "完成"

# ------------------------------------------------------------------
# 2. Add the new weekly block (rows 172-180), cloning the layout/
#    formatting of the preceding block (rows 162-170) which has the
#    exact same structure (header / column-titles / 6 data rows /
#    summary row) and merged cells.
# ------------------------------------------------------------------
$src = $ws.Range("A162:D170")
$dst = $ws.Range("A172:D180")
$src.Copy($dst)

# Header / date row
$ws.Range("A172").Value = "日期：2017.10.30 第十周 周一"

# Data rows - update the "计划任务" (plan) column with the new entries
$ws.Range("B174").Value = "继续开发后台“订单管理”模块中的实体类和控制层(controller)"
$ws.Range("B175").Value = "完成首页“分类”模块以及开发前端首页的“宝贝分类”模块"
$ws.Range("B176").Value = "帮助前端开发人员开发其中一个小模块"
$ws.Range("B177").Value = "继续开发前端“我的”模块"
$ws.Range("B178").Value = "继续开发后台“宝贝管理”模块中的控制层(controller)"
$ws.Range("B179").Value = "帮助前端开发人员开发其中一个小模块"

# "完成情况" column for the new block is left blank, so clear whatever got
# copied over from the source block.
$ws.Range("C174").ClearContents()
$ws.Range("C175").ClearContents()
$ws.Range("C176").ClearContents()
$ws.Range("C177").ClearContents()
$ws.Range("C178").ClearContents()
$ws.Range("C179").ClearContents()

# Summary row
$ws.Range("A180").Value = "总结："

# Row heights for the newly added rows
$ws.Rows.Item(172).RowHeight = 22.5
$ws.Rows.Item(173).RowHeight = 22.5
$ws.Rows.Item(174).RowHeight = 67.5
$ws.Rows.Item(175).RowHeight = 67.5
$ws.Rows.Item(176).RowHeight = 45
$ws.Rows.Item(177).RowHeight = 22.5
$ws.Rows.Item(178).RowHeight = 67.5
$ws.Rows.Item(179).RowHeight = 45
$ws.Rows.Item(180).RowHeight = 22.5

# ------------------------------------------------------------------
# 3. Update the sheet view - scroll position & selected cell
# ------------------------------------------------------------------
$ws.Range("B179").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 168
$win.ScrollColumn = 1
